# The commit adds one new weekly price-report record for "Albahaca" at the
# Vega Central Mapocho de Santiago market. In the underlying data the rows
# are kept in (reverse-ish) date order, so the new record is inserted as
# row 534, pushing the former rows 534-651 down to 535-652 and growing the
# sheet's used range from A1:R651 to A1:R652.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 534..651 down by one to make room for the new record.
$ws.Rows.Item(534).Insert()

# Populate the newly inserted row 534 with the new record's data.
$ws.Cells.Item(534, 1).Value  = 9
$ws.Cells.Item(534, 2).Value  = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(534, 3).Value  = 'Metropolitana'
$ws.Cells.Item(534, 4).Value  = 45244
$ws.Cells.Item(534, 5).Value  = 13
$ws.Cells.Item(534, 6).Value  = 100112052
$ws.Cells.Item(534, 7).Value  = 'Albahaca'
$ws.Cells.Item(534, 8).Value  = 'Sin especificar'
$ws.Cells.Item(534, 9).Value  = 'Primera'
$ws.Cells.Item(534, 10).Value = 160
$ws.Cells.Item(534, 11).Value = 5000
$ws.Cells.Item(534, 12).Value = 6000
$ws.Cells.Item(534, 13).Value = 5500
$ws.Cells.Item(534, 14).Value = '$/docena de matas'
$ws.Cells.Item(534, 15).Value = 'Provincia de Chacabuco'
$ws.Cells.Item(534, 16).Value = 917
$ws.Cells.Item(534, 17).Value = 6
$ws.Cells.Item(534, 18).Value = 'Hortaliza'
